# Auto-generated edit script applying numeric cell updates to the Leve profit
# tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching the
# refreshed currentAveragePrice / LevePrice / LeveProfit figures from the
# scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 301.84616
$ws.Range("I2").Value = 546.8333
$ws.Range("J2").Value = 91.85714
$ws.Range("K2").Value = 546.8333
$ws.Range("L2").Value = 91.85714
$ws.Range("M2").Value = -433.8333
$ws.Range("N2").Value = -317.85714
$ws.Range("H4").Value = 74.75
$ws.Range("I4").Value = 66.333336
$ws.Range("K4").Value = 66.333336
$ws.Range("M4").Value = 47.666664
$ws.Range("H5").Value = 99.5
$ws.Range("I5").Value = 99.5
$ws.Range("K5").Value = 99.5
$ws.Range("M5").Value = 15.5
$ws.Range("H17").Value = 2100
$ws.Range("J17").Value = 2100
$ws.Range("L17").Value = 6300
$ws.Range("N17").Value = -6636
$ws.Range("H32").Value = 1499.5
$ws.Range("H33").Value = 96.8125
$ws.Range("I33").Value = 72.63636
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 72.63636
$ws.Range("L33").Value = 150
$ws.Range("M33").Value = 156.36364
$ws.Range("N33").Value = -608
$ws.Range("H69").Value = 8490
$ws.Range("J69").Value = 8490
$ws.Range("L69").Value = 25470
$ws.Range("N69").Value = -27218
$ws.Range("H72").Value = 8490
$ws.Range("J72").Value = 8490
$ws.Range("L72").Value = 76410
$ws.Range("N72").Value = -85146
$ws.Range("H125").Value = 1998
$ws.Range("I125").Value = 1497
$ws.Range("K125").Value = 13473
$ws.Range("M125").Value = -11013
$ws.Range("H137").Value = 3725.5
$ws.Range("J137").Value = 5000
$ws.Range("L137").Value = 15000
$ws.Range("N137").Value = -20100
$ws.Range("H138").Value = 3243.4285
$ws.Range("I138").Value = 2891.1667
$ws.Range("J138").Value = 3384.3333
$ws.Range("K138").Value = 8673.500100000001
$ws.Range("L138").Value = 10152.9999
$ws.Range("M138").Value = -3533.500100000001
$ws.Range("N138").Value = -20432.9999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 100.4
$ws.Range("I5").Value = 88
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 88
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = 24
$ws.Range("N5").Value = -374
$ws.Range("H104").Value = 20806.25
$ws.Range("J104").Value = 20806.25
$ws.Range("L104").Value = 20806.25
$ws.Range("N104").Value = -27794.25
$ws.Range("H118").Value = 50000
$ws.Range("J118").Value = 50000
$ws.Range("L118").Value = 50000
$ws.Range("N118").Value = -53314
$ws.Range("H122").Value = 998.1539
$ws.Range("I122").Value = 995.2
$ws.Range("J122").Value = 1008
$ws.Range("K122").Value = 2985.6
$ws.Range("L122").Value = 3024
$ws.Range("M122").Value = -535.6000000000004
$ws.Range("N122").Value = -7924
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 100.4
$ws.Range("I4").Value = 88
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 88
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = 27
$ws.Range("N4").Value = -380
$ws.Range("H10").Value = 1796
$ws.Range("I10").Value = 2199
$ws.Range("J10").Value = 990
$ws.Range("K10").Value = 2199
$ws.Range("L10").Value = 990
$ws.Range("M10").Value = -2059
$ws.Range("N10").Value = -1270
$ws.Range("H54").Value = 4846.3335
$ws.Range("I54").Value = 4846.3335
$ws.Range("K54").Value = 4846.3335
$ws.Range("M54").Value = -4362.3335
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 40.6
$ws.Range("I7").Value = 44.333332
$ws.Range("K7").Value = 44.333332
$ws.Range("M7").Value = 68.666668
$ws.Range("H8").Value = 1245
$ws.Range("I8").Value = 725
$ws.Range("J8").Value = 1375
$ws.Range("K8").Value = 725
$ws.Range("L8").Value = 1375
$ws.Range("M8").Value = -585
$ws.Range("N8").Value = -1655
$ws.Range("H32").Value = 1459.1818
$ws.Range("I32").Value = 479.42856
$ws.Range("J32").Value = 3173.75
$ws.Range("K32").Value = 479.42856
$ws.Range("L32").Value = 3173.75
$ws.Range("M32").Value = -163.42856
$ws.Range("N32").Value = -3805.75
$ws.Range("H35").Value = 1750
$ws.Range("I35").Value = 1125.5
$ws.Range("J35").Value = 2999
$ws.Range("K35").Value = 1125.5
$ws.Range("L35").Value = 2999
$ws.Range("M35").Value = -831.5
$ws.Range("N35").Value = -3587
$ws.Range("H42").Value = 11166
$ws.Range("I42").Value = 4500
$ws.Range("J42").Value = 14499
$ws.Range("K42").Value = 4500
$ws.Range("L42").Value = 14499
$ws.Range("M42").Value = -3907
$ws.Range("N42").Value = -15685
$ws.Range("H141").Value = 855971.25
$ws.Range("J141").Value = 1121295.4
$ws.Range("L141").Value = 1121295.4
$ws.Range("N141").Value = -1131655.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 974.8889
$ws.Range("I2").Value = 389.86667
$ws.Range("K2").Value = 2339.20002
$ws.Range("M2").Value = -2226.20002
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2888
$ws.Range("H11").Value = 3738.889
$ws.Range("I11").Value = 300
$ws.Range("J11").Value = 4168.75
$ws.Range("K11").Value = 900
$ws.Range("L11").Value = 12506.25
$ws.Range("M11").Value = -760
$ws.Range("N11").Value = -12786.25
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = $null
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = $null
$ws.Range("H99").Value = 766.6667
$ws.Range("I99").Value = 750
$ws.Range("J99").Value = 775
$ws.Range("K99").Value = 2250
$ws.Range("L99").Value = 2325
$ws.Range("M99").Value = -4
$ws.Range("N99").Value = -6817
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 306.33334
$ws.Range("I13").Value = 112.666664
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 112.666664
$ws.Range("L13").Value = 500
$ws.Range("M13").Value = 26.333336
$ws.Range("N13").Value = -778
$ws.Range("H22").Value = 1365.8334
$ws.Range("I22").Value = 733.3333
$ws.Range("J22").Value = 1998.3334
$ws.Range("K22").Value = 733.3333
$ws.Range("L22").Value = 1998.3334
$ws.Range("M22").Value = -204.3333
$ws.Range("N22").Value = -3056.3334
$ws.Range("H23").Value = 3150
$ws.Range("J23").Value = 3150
$ws.Range("L23").Value = 3150
$ws.Range("N23").Value = -3596
$ws.Range("H93").Value = 72495
$ws.Range("J93").Value = 70000
$ws.Range("L93").Value = 70000
$ws.Range("N93").Value = -73744
$ws.Range("H102").Value = 4211.8887
$ws.Range("J102").Value = 3500
$ws.Range("L102").Value = 3500
$ws.Range("N102").Value = -6744
$ws.Range("H105").Value = 16750
$ws.Range("J105").Value = 16750
$ws.Range("L105").Value = 16750
$ws.Range("N105").Value = -23738
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 30000
$ws.Range("J43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30386
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = $null
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = $null
$ws.Range("H100").Value = 2998
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").Value = $null
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null
$ws.Range("H132").Value = 10000
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = $null
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 6924.8335
$ws.Range("I4").Value = 5312.25
$ws.Range("J4").Value = 10150
$ws.Range("K4").Value = 5312.25
$ws.Range("L4").Value = 10150
$ws.Range("M4").Value = -5199.25
$ws.Range("N4").Value = -10376
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = $null
$ws.Range("H70").Value = 21500
$ws.Range("I70").Value = 18000
$ws.Range("K70").Value = 18000
$ws.Range("M70").Value = -17685
$ws.Range("H73").Value = 21500
$ws.Range("I73").Value = 18000
$ws.Range("K73").Value = 18000
$ws.Range("M73").Value = -16908
$ws.Range("H75").Value = 90000
$ws.Range("I75").Value = 90000
$ws.Range("K75").Value = 90000
$ws.Range("M75").Value = -89064
$ws.Range("H78").Value = 90000
$ws.Range("I78").Value = 90000
$ws.Range("K78").Value = 270000
$ws.Range("M78").Value = -265320
$ws.Range("H107").Value = 1395.8334
$ws.Range("J107").Value = 1422.1428
$ws.Range("L107").Value = 4266.428400000001
$ws.Range("N107").Value = -8106.428400000001
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null
